$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header: I1 = "minority_percent" ---
$ws.Range("I1").Value = "minority_percent"

# Give I1 the same bold/boxed header look already used by A1:H1 (style 1),
# then trim the border down to left/right only (matching the new style added
# to styles.xml for this column's header).
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$cell = $ws.Range("I1")
$cell.Borders.Item(8).LineStyle = -4142   # xlEdgeTop    -> none
$cell.Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none
$cell.Interior.Pattern = -4142            # xlNone (explicit "no fill")

# --- New column data: minority_percent per dataset group ---
$values = @(0.3,0.3,0.3,0.3,0.3,0.3,0.2,0.2,0.2,0.2,0.2,0.2,0.1,0.1,0.1,0.1,0.1,0.1,0.05,0.05,0.05,0.05,0.05,0.05,0.01,0.01,0.01,0.01,0.01,0.01,0.005,0.005,0.005,0.005,0.005,0.005)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i]
}

# Column I needs to be wide enough for its header/content, like the other cols.
$ws.Columns.Item(9).AutoFit()

# Scroll back to the top and select the newly added data range, matching the
# saved view state.
$ws.Range("I32:I37").Select() | Out-Null
